# "fixup! Update export fixture files"
#
# Three runs in the fixture document have a stray literal newline baked
# into their <w:t> text (an artifact of how the fixture was generated).
# This normalizes them:
#   - "This is an annotatable resource in the casebook.\n"   -> no trailing \n
#   - "This is the second chapter of the casebook.\n"        -> no trailing \n
#   - the long "highlighted: ... ;" paragraph has two internal \n line
#     breaks that should become plain spaces so it reads as one line.

$d = $word.ActiveDocument

function Replace-InParagraph($paraIndex, $findText, $replaceText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, `
        $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        throw "Replace-InParagraph: text not found in paragraph $paraIndex : $findText"
    }
}

# --- "This is an annotatable resource in the casebook.\n" -------------
Replace-InParagraph 23 "casebook.`n" "casebook."

# --- "highlighted: ...elide;\nreplaced: ... noted:\ncontent to note; ..." --
Replace-InParagraph 24 "elide;`nreplaced" "elide; replaced"
Replace-InParagraph 24 "noted:`ncontent"  "noted: content"

# --- "This is the second chapter of the casebook.\n" -------------------
Replace-InParagraph 28 "casebook.`n" "casebook."
